$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Insert a new row above row 28 (shifts the Science 37 row, currently row 28, down to row 29)
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the Avantor entry (plain text, no style/hyperlink, matching sibling rows 26/27)
$ws.Cells.Item(28, 2).Value = "Avantor"
$ws.Cells.Item(28, 3).Value = "AVTR"

# The row insert above does not relocate the worksheet's hyperlink anchors, so
# rebuild the hyperlink collection to point at the shifted cells.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B29"), "SNCE.xlsx") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "UNH.xlsx") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "CVS.xlsx") | Out-Null

# Hyperlinks.Add mints a fresh "Hyperlink" style entry instead of reusing the
# workbook's existing one; re-apply the named style so these cells keep
# pointing at the original shared style index.
$ws.Range("B29").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"

# Update selection to match the post-edit active cell shown in the diff
$ws.Range("D28").Select() | Out-Null
